# Insert a new data row at sheet row 268 (pushing existing rows 268-381 down to 269-382)
# and populate it with the new Brócoli price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

$ws.Cells.Item(268, 1).Value = 4
$ws.Cells.Item(268, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(268, 3).Value = "Los Lagos"
$ws.Cells.Item(268, 4).Value = 44784
$ws.Cells.Item(268, 5).Value = 10
$ws.Cells.Item(268, 6).Value = 100112023
$ws.Cells.Item(268, 7).Value = "Brócoli"
$ws.Cells.Item(268, 8).Value = "Sin especificar"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 500
$ws.Cells.Item(268, 11).Value = 1600
$ws.Cells.Item(268, 12).Value = 1600
$ws.Cells.Item(268, 13).Value = 1600
$ws.Cells.Item(268, 14).Value = "`$/unidad"
$ws.Cells.Item(268, 15).Value = "Región Metropolitana"
$ws.Cells.Item(268, 16).Value = 1600
$ws.Cells.Item(268, 17).Value = 1
$ws.Cells.Item(268, 18).Value = "Hortaliza"
